# Handle duplicate sample IDs & values
# Update the funding/organization name on the Personnel sheet for the
# OOI CGSN Data Team row (D6) from "Ocean Observatories Initiative"
# to "NSF Ocean Observatories Initiative".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

$ws.Range("D6").Value = "NSF Ocean Observatories Initiative"

$ws.Activate()
$ws.Range("D16").Select()
